# This edit inserts two new daily price records (rows) for "Mango" at
# Vega Central Mapocho de Santiago, right before the existing row that used
# to be row 220 (date 2020-12-07 / serial 44172). All rows from the old
# row 220 through the old last row (295) shift down by two positions
# (to new rows 222-297), and the sheet's used range grows from A1:T295
# to A1:T297.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 220 (pushes existing row 220 and
# everything below it down by two rows).
$ws.Rows.Item(220).Insert()
$ws.Rows.Item(220).Insert()

# ---- New row 220 ----
$ws.Cells.Item(220,1).Value  = 9
$ws.Cells.Item(220,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(220,3).Value  = "Metropolitana"
$ws.Cells.Item(220,4).Value  = 44524
$ws.Cells.Item(220,5).Value  = 13
$ws.Cells.Item(220,6).Value  = "Fruta"
$ws.Cells.Item(220,7).Value  = 100108
$ws.Cells.Item(220,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(220,9).Value  = 100108002
$ws.Cells.Item(220,10).Value = "Mango"
$ws.Cells.Item(220,11).Value = "Sin especificar"
$ws.Cells.Item(220,12).Value = "Primera"
$ws.Cells.Item(220,13).Value = 300
$ws.Cells.Item(220,14).Value = 6000
$ws.Cells.Item(220,15).Value = 6000
$ws.Cells.Item(220,16).Value = 6000
$ws.Cells.Item(220,17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(220,18).Value = "Brasil"
$ws.Cells.Item(220,19).Value = 1500
$ws.Cells.Item(220,20).Value = 4

# ---- New row 221 ----
$ws.Cells.Item(221,1).Value  = 9
$ws.Cells.Item(221,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(221,3).Value  = "Metropolitana"
$ws.Cells.Item(221,4).Value  = 44524
$ws.Cells.Item(221,5).Value  = 13
$ws.Cells.Item(221,6).Value  = "Fruta"
$ws.Cells.Item(221,7).Value  = 100108
$ws.Cells.Item(221,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(221,9).Value  = 100108002
$ws.Cells.Item(221,10).Value = "Mango"
$ws.Cells.Item(221,11).Value = "Sin especificar"
$ws.Cells.Item(221,12).Value = "Primera"
$ws.Cells.Item(221,13).Value = 530
$ws.Cells.Item(221,14).Value = 6000
$ws.Cells.Item(221,15).Value = 6500
$ws.Cells.Item(221,16).Value = 6236
$ws.Cells.Item(221,17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(221,18).Value = "Perú"
$ws.Cells.Item(221,19).Value = 1559
$ws.Cells.Item(221,20).Value = 4
